# aggiornamento a l 23 agosto 2021
# Appends new daily rows (344-357, 2021-08-10 .. 2021-08-23) to the
# COVID-style tracking sheet: date, nuovi pos., somma mobile 7gg.,
# somma mobile 7gg. per 100mila abitanti.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# date serial, nuovi pos. (B), somma mobile 7gg. (C), per 100mila (D)
$newRows = @(
    @(44418, 1, 1, 31.25976867771178),
    @(44419, 0, 1, 31.25976867771178),
    @(44420, 0, 1, 31.25976867771178),
    @(44421, 0, 1, 31.25976867771178),
    @(44422, 1, 2, 62.51953735542357),
    @(44423, 1, 3, 93.77930603313536),
    @(44424, 1, 4, 125.0390747108471),
    @(44425, 0, 3, 93.77930603313536),
    @(44426, 0, 3, 93.77930603313536),
    @(44427, 0, 3, 93.77930603313536),
    @(44428, 1, 4, 125.0390747108471),
    @(44429, 0, 3, 93.77930603313536),
    @(44430, 0, 2, 62.51953735542357),
    @(44431, 0, 1, 31.25976867771178)
)

$lastExistingRow = 343
$startRow = $lastExistingRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

# Column A carries the date style/number format used by the rest of
# the table; replicate it onto the freshly added date cells.
$endRow = $startRow + $newRows.Count - 1
$styleSource = $ws.Range("A" + $lastExistingRow)
$styleTarget = $ws.Range("A" + $startRow + ":A" + $endRow)
$styleSource.Copy()
$styleTarget.PasteSpecial(-4122)
